# Apply odds updates to Sheet1 as described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("F2").Value = 1.57
$ws.Range("K2").Value = 5
$ws.Range("Q2").Value = 1.78

# Row 3
$ws.Range("H3").Value = 1.41
$ws.Range("Q3").Value = 1.58
$ws.Range("AI3").Value = 34

# Row 4
$ws.Range("F4").Value = 1.3
$ws.Range("G4").Value = 1.71
$ws.Range("H4").Value = 1.16
$ws.Range("I4").Value = 44
$ws.Range("J4").Value = 1.36
$ws.Range("K4").Value = 32
$ws.Range("Q4").Value = 1.66
$ws.Range("V4").Value = 1.05

# Row 5
$ws.Range("H5").Value = 3.15
$ws.Range("I5").Value = 4.6
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 5.8
$ws.Range("V5").Value = 1.3
$ws.Range("W5").Value = 1.63

# Row 6
$ws.Range("F6").Value = 1.77
$ws.Range("G6").Value = 1.94
$ws.Range("J6").Value = 3.3
$ws.Range("K6").Value = 3.75
$ws.Range("M6").Value = 1.09
$ws.Range("P6").Value = 1.69
$ws.Range("T6").Value = 2
$ws.Range("V6").Value = 1.2
$ws.Range("W6").Value = 2.06
$ws.Range("X6").Value = 13.5
$ws.Range("Z6").Value = 55
$ws.Range("AJ6").Value = 25
